$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 text content ---
$ws.Range("B2").Value = "Sub Category "
$ws.Range("C2").Value = 'COLD & FLU ; IBS ; KIDS COLD & FLU ; COLD & FLU ;KIDS COUGH ; KIDS DECONGESTANTS ; KIDS DIGESTIVE HEALTH ; KIDS HAYFEVER ; LAXATIVES ; KIDS TEETHING ; WIND ; ALLERGY ; REHYDRATION ; PROBIOTICS ; PAIN MANAGEMENT ; ANTI-AGE FACE ; ARTIFICIAL TAN ; BODY CLEANSING ; BODY SPRAY ; COSMETICS ; DEODORANTS ; MEN''S TOILETRIES ; SUNCARE ; HAIR CARE ;  BABY HEALTHCARE ; FOR MUM ; KIDS HAIRCARE ; KIDS TOILETRIES ;  KIDS WIPES ; COTTON ; INCONTINENCE ; SANITARY TOWELS ; FEMININE WASH ; KIDS MOUTHWASH '

# --- Remove now-unused trailing exclusion columns (location_type / template_name / brand_name data) ---
$ws.Range("D2:I2").ClearContents()

# --- Delete the Exclude3/Value3/Exclude4/Value4 header columns (F:I) ---
$ws.Range("F1:I2").EntireColumn.Delete()

# --- Alignment updates for row 2 ---
$ws.Range("A2:B2").VerticalAlignment = -4108  # xlCenter
$ws.Range("C2").VerticalAlignment = -4107     # xlBottom (default)
$ws.Range("C2").WrapText = $true

# --- Row height ---
$ws.Range("A2").RowHeight = 247.4

# --- Column widths (character units); engine quantizes to Excel's internal pixel grid ---
$ws.Range("A1").ColumnWidth = 32.166666666666664
$ws.Range("B1").ColumnWidth = 15.999999999999998
$ws.Range("C1").ColumnWidth = 24.5
$ws.Range("D1").ColumnWidth = 15.166666666666666
$ws.Range("E1").ColumnWidth = 23.833333333333336

# --- Selection ---
$null = $ws.Range("A8").Select()

# --- Cosmetic tab-ratio tweak (window split between sheet tabs / scrollbar) ---
$excel.ActiveWindow.TabRatio = 0.993
